$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two obsolete sub-sector rows (18_01_02_01_gasturbine, 18_01_02_02_combinedcycle)
$ws.Rows("247:248").Delete()

# The row that used to be 249 (18_01_02_03_ccs) is now row 247; rename it to the new sub-sector
$ws.Range("B247").Value = "18_01_02_gas_power_ccs"
$ws.Range("C247").Value = "18_01_02_gas_power_ccs"
